$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1195
$ws.Range("I19").Value = 1017.2
$ws.Range("K19").Value = 1017.2
$ws.Range("M19").Value = -842.2
$ws.Range("H43").Value = 6295.2
$ws.Range("I43").Value = 5501
$ws.Range("J43").Value = 6493.75
$ws.Range("K43").Value = 5501
$ws.Range("L43").Value = 6493.75
$ws.Range("M43").Value = -5432
$ws.Range("N43").Value = -6631.75
$ws.Range("H98").Value = 997.5
$ws.Range("I98").Value = 997
$ws.Range("K98").Value = 997
$ws.Range("M98").Value = 501
$ws.Range("H116").Value = 4693.5713
$ws.Range("I116").Value = 4772.125
$ws.Range("J116").Value = 4588.8335
$ws.Range("K116").Value = 4772.125
$ws.Range("L116").Value = 4588.8335
$ws.Range("M116").Value = -1330.125
$ws.Range("N116").Value = -11472.8335
$ws.Range("H122").Value = 997.5
$ws.Range("I122").Value = 997
$ws.Range("K122").Value = 2991
$ws.Range("M122").Value = -541
$ws.Range("H132").Value = 7798
$ws.Range("I132").Value = 7961.4546
$ws.Range("K132").Value = 23884.3638
$ws.Range("M132").Value = -21354.3638
$ws.Range("H138").Value = 2741.2778
$ws.Range("J138").Value = 2981.6428
$ws.Range("L138").Value = 8944.928400000001
$ws.Range("N138").Value = -19224.9284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1124.2
$ws.Range("I2").Value = 1252
$ws.Range("K2").Value = 1252
$ws.Range("M2").Value = -1139
$ws.Range("H43").Value = 100377
$ws.Range("J43").Value = 100377
$ws.Range("L43").Value = 100377
$ws.Range("N43").Value = -101003
$ws.Range("H45").Value = 1011
$ws.Range("I45").Value = 1011
$ws.Range("K45").Value = 1011
$ws.Range("M45").Value = -634
$ws.Range("H61").Value = 2209.6924
$ws.Range("I61").Value = 2209.6924
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2209.6924
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1997.6924
$ws.Range("N61").ClearContents()
$ws.Range("H116").Value = 1124.2
$ws.Range("I116").Value = 1252
$ws.Range("K116").Value = 1252
$ws.Range("M116").Value = 1042
$ws.Range("H122").Value = 1571.125
$ws.Range("I122").Value = 1571.125
$ws.Range("K122").Value = 4713.375
$ws.Range("M122").Value = -2263.375
$ws.Range("H132").Value = 4268.5
$ws.Range("I132").Value = 4011.2
$ws.Range("K132").Value = 12033.6
$ws.Range("M132").Value = -9503.599999999999
$ws.Range("H136").Value = 2209.6924
$ws.Range("I136").Value = 2209.6924
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6629.0772
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4079.0772
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1124.2
$ws.Range("I3").Value = 1252
$ws.Range("K3").Value = 1252
$ws.Range("M3").Value = -1138
$ws.Range("H99").Value = 2348.35
$ws.Range("I99").Value = 1531.4445
$ws.Range("K99").Value = 1531.4445
$ws.Range("M99").Value = -33.44450000000006
$ws.Range("H107").Value = 1764.1765
$ws.Range("I107").Value = 1764.1765
$ws.Range("K107").Value = 1764.1765
$ws.Range("M107").Value = 155.8235
$ws.Range("H134").Value = 1111.1428
$ws.Range("I134").Value = 1154.9231
$ws.Range("K134").Value = 3464.7693
$ws.Range("M134").Value = -929.7692999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4398.5386
$ws.Range("I86").Value = 4390.0835
$ws.Range("K86").Value = 4390.0835
$ws.Range("M86").Value = -3267.0835
$ws.Range("H89").Value = 4398.5386
$ws.Range("I89").Value = 4390.0835
$ws.Range("K89").Value = 21950.4175
$ws.Range("M89").Value = -16334.4175
$ws.Range("H122").Value = 4160.222
$ws.Range("I122").Value = 4293
$ws.Range("J122").Value = 3994.25
$ws.Range("K122").Value = 12879
$ws.Range("L122").Value = 11982.75
$ws.Range("M122").Value = -10429
$ws.Range("N122").Value = -16882.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 104.166664
$ws.Range("I23").Value = 87.22221999999999
$ws.Range("J23").Value = 155
$ws.Range("K23").Value = 261.66666
$ws.Range("L23").Value = 465
$ws.Range("M23").Value = -26.66665999999998
$ws.Range("N23").Value = -935
$ws.Range("H34").Value = 4888.4614
$ws.Range("J34").Value = 5691
$ws.Range("L34").Value = 17073
$ws.Range("N34").Value = -17241
$ws.Range("H39").Value = 7800
$ws.Range("J39").Value = 7800
$ws.Range("L39").Value = 23400
$ws.Range("N39").Value = -23988
$ws.Range("H50").Value = 3118.25
$ws.Range("I50").Value = 265
$ws.Range("K50").Value = 795
$ws.Range("M50").Value = -314
$ws.Range("H53").Value = 3118.25
$ws.Range("I53").Value = 265
$ws.Range("K53").Value = 795
$ws.Range("M53").Value = -314
$ws.Range("H55").Value = 2575
$ws.Range("I55").Value = 1650
$ws.Range("J55").Value = 3500
$ws.Range("K55").Value = 4950
$ws.Range("L55").Value = 10500
$ws.Range("M55").Value = -4773
$ws.Range("N55").Value = -10854
$ws.Range("H61").Value = 132.33333
$ws.Range("I61").Value = 152.83333
$ws.Range("J61").Value = 91.333336
$ws.Range("K61").Value = 458.49999
$ws.Range("L61").Value = 274.000008
$ws.Range("M61").Value = -243.49999
$ws.Range("N61").Value = -704.000008
$ws.Range("H94").Value = 7924.625
$ws.Range("J94").Value = 7924.625
$ws.Range("L94").Value = 23773.875
$ws.Range("N94").Value = -25125.875
$ws.Range("H128").Value = 394999.5
$ws.Range("I128").Value = 394999.5
$ws.Range("K128").Value = 1184998.5
$ws.Range("M128").Value = -1180018.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 47474.75
$ws.Range("I63").Value = 50000
$ws.Range("J63").Value = 44949.5
$ws.Range("K63").Value = 50000
$ws.Range("L63").Value = 44949.5
$ws.Range("M63").Value = -49314
$ws.Range("N63").Value = -46321.5
$ws.Range("H66").Value = 47474.75
$ws.Range("I66").Value = 50000
$ws.Range("J66").Value = 44949.5
$ws.Range("K66").Value = 150000
$ws.Range("L66").Value = 134848.5
$ws.Range("M66").Value = -146568
$ws.Range("N66").Value = -141712.5
$ws.Range("H97").Value = 325.78946
$ws.Range("I97").Value = 303.75
$ws.Range("J97").Value = 443.33334
$ws.Range("K97").Value = 303.75
$ws.Range("L97").Value = 443.33334
$ws.Range("M97").Value = 192.25
$ws.Range("N97").Value = -1435.33334
$ws.Range("H107").Value = 2152.4546
$ws.Range("I107").Value = 717.8
$ws.Range("J107").Value = 3348
$ws.Range("K107").Value = 717.8
$ws.Range("L107").Value = 3348
$ws.Range("M107").Value = 1202.2
$ws.Range("N107").Value = -7188
$ws.Range("H122").Value = 1613.5555
$ws.Range("I122").Value = 1613.5555
$ws.Range("K122").Value = 4840.666499999999
$ws.Range("M122").Value = -2390.666499999999
$ws.Range("H132").Value = 1921
$ws.Range("I132").Value = 1824.5
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5473.5
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -2943.5
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4610.8
$ws.Range("I46").Value = 3650
$ws.Range("K46").Value = 3650
$ws.Range("M46").Value = -3462
$ws.Range("H82").Value = 4195.9
$ws.Range("J82").Value = 7174.75
$ws.Range("L82").Value = 7174.75
$ws.Range("N82").Value = -7896.75
$ws.Range("H85").Value = 4195.9
$ws.Range("J85").Value = 7174.75
$ws.Range("L85").Value = 7174.75
$ws.Range("N85").Value = -9670.75
$ws.Range("H100").Value = 9750
$ws.Range("I100").Value = 9500
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 9500
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -8959
$ws.Range("N100").Value = -11082
$ws.Range("H132").Value = 3057.1428
$ws.Range("I132").Value = 3033.5
$ws.Range("K132").Value = 9100.5
$ws.Range("M132").Value = -6570.5
$ws.Range("H136").Value = 71504856
$ws.Range("I136").Value = 63499
$ws.Range("J136").Value = 166760000
$ws.Range("K136").Value = 190497
$ws.Range("L136").Value = 500280000
$ws.Range("M136").Value = -187947
$ws.Range("N136").Value = -500285100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3168.5715
$ws.Range("I132").Value = 3280
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 9840
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -7310
$ws.Range("N132").Value = -12560
